# HBVCALIBandGages.xlsx - add the OakGrove40 HBVCALIB / gage record to the
# "Clackamas" sheet as a primary listed row (row 9), promoted out of the
# "Other USGS gages" sub-list (which previously started at row 11). The old
# row 11 record is removed (its data now lives in row 9), so the records
# that used to be rows 12-14 shift up to rows 11-13.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Clackamas")
$ws.Activate()

# Preserve the "Downloaded Flow/Temp" column formatting (centered custom
# text format) from the old row 11 before it is removed, and apply it to
# the new row 9 so the promoted record keeps the same look as its peers.
$ws.Range("P11:Q11").Copy() | Out-Null
$ws.Range("P9:Q9").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = $false

# Remove the old "OAK GROVE FORK NEAR GOVERNMENT CAMP, OR" row from the
# "Other USGS gages" list -- its data is about to be re-entered above as a
# primary HBVCALIB/gage row. Deleting shifts the rows below it (old 12-14)
# up by one, to 11-13.
$ws.Rows("11:11").Delete()

# Re-enter that gage's data as a new, primary row 9, with the extra
# HBVCALIB identity columns (A/B/C) that the other primary rows (7, 8)
# have.
$ws.Range("A9").Value = 40
$ws.Range("B9").Value = "OakGrove40"
$ws.Range("C9").Value = 23809450
$ws.Range("H9").Value = 14208700
$ws.Range("I9").Value = 14143
$ws.Range("J9").Value = 0
$ws.Range("L9").Value = " OAK GROVE FORK NEAR GOVERNMENT CAMP, OR"
$ws.Range("M9").Value = 23810706
$ws.Range("N9").Value = 54.4
$ws.Range("P9").Value = "X"
$ws.Range("Q9").Value = "M"

# Match the sheet's new active selection.
$ws.Range("D9").Select()
